$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header): add P1 / Q1, matching style of the existing O1 header cell ---
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("P1").Value = 14
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("Q1").Value = 15

# --- Rows 2-25: swap I/K/M/O values and append P/Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # column I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # column K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # column M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # column O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # column P (new)
    $ws.Cells.Item($r, 17).Value = 2  # column Q (new)
}
